# Locate the paragraph that ends with "Meeting: October 27th, 2015 in MSB 6503"
# and add a new BodyText paragraph right after it, containing the new
# introductory sentences drafted for the "big picture" intro.

$d = $word.ActiveDocument

$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Meeting: October 27th, 2015 in MSB 6503*") {
        $targetPara = $p
        break
    }
}

# Insert a brand-new paragraph immediately after the target paragraph.
$targetPara.Range.InsertParagraphAfter()

# The freshly-created paragraph is the next one in the document.
$newPara = $targetPara.Next()

$newPara.Range.Text = "I am interested in understanding how bacterial populations diversify and form species. My work aims to better understand the forces that shape bacterial population structure by using the lakes as a model system. The McMahon Lab metagenomic time series of many lakes is an ideal dataset for investigating how bacterial populations change through time."
$newPara.Style = "BodyText"
